# Horarios actualizados Linea 141 - 1221
# Refresh of the scraped schedule data across the three worksheets
# (LP1912, LP1912-215, 6203-6173). A new scrape ran at 05:23:04, so the
# "Ultima actualizacion" / "Total filas" summary cells and the data rows
# are rewritten to match the new snapshot (rows shift down as older
# departures fall off / new ones are appended).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912" ---------------------------------------------------
    $ws1.Cells.Item(2,1).Value = "Última actualización: 05:23:04"

    $ws1.Cells.Item(3,1).Value = "Total filas: 30"

    $ws1.Cells.Item(10,1).Value = "05:23:04"
    $ws1.Cells.Item(10,2).Value = "05:23"
    $ws1.Cells.Item(10,3).Value = "23_HERNANDEZ"
    $ws1.Cells.Item(10,4).Value = 0

    $ws1.Cells.Item(11,1).Value = "05:23:04"
    $ws1.Cells.Item(11,2).Value = "05:32"
    $ws1.Cells.Item(11,3).Value = "81_EL PELIGRO"
    $ws1.Cells.Item(11,4).Value = 9

    $ws1.Cells.Item(12,2).Value = "05:34"
    $ws1.Cells.Item(12,3).Value = "215B_EL PATO"
    $ws1.Cells.Item(12,4).Value = 38

    $ws1.Cells.Item(13,1).Value = "05:23:04"
    $ws1.Cells.Item(13,2).Value = "05:44"
    $ws1.Cells.Item(13,3).Value = "14_ABASTO"
    $ws1.Cells.Item(13,4).Value = 21

    $ws1.Cells.Item(14,2).Value = "05:46"
    $ws1.Cells.Item(14,3).Value = "15_ABASTO"
    $ws1.Cells.Item(14,4).Value = 50

    $ws1.Cells.Item(15,1).Value = "05:23:04"
    $ws1.Cells.Item(15,2).Value = "05:52"
    $ws1.Cells.Item(15,3).Value = "17_ROMERO"
    $ws1.Cells.Item(15,4).Value = 29

    $ws1.Cells.Item(16,2).Value = "05:54"
    $ws1.Cells.Item(16,3).Value = "10_OLMOS"
    $ws1.Cells.Item(16,4).Value = 58

    $ws1.Cells.Item(17,1).Value = "05:23:04"
    $ws1.Cells.Item(17,2).Value = "06:04"
    $ws1.Cells.Item(17,3).Value = "10_OLMOS"
    $ws1.Cells.Item(17,4).Value = 41

    $ws1.Cells.Item(18,2).Value = "06:04"
    $ws1.Cells.Item(18,3).Value = "16_SANTA ANA"
    $ws1.Cells.Item(18,4).Value = 68

    $ws1.Cells.Item(19,1).Value = "05:23:04"
    $ws1.Cells.Item(19,2).Value = "06:11"
    $ws1.Cells.Item(19,3).Value = "215A_EL PATO"
    $ws1.Cells.Item(19,4).Value = 48

    $ws1.Cells.Item(20,2).Value = "06:14"
    $ws1.Cells.Item(20,3).Value = "225_HARAS DEL SUR"
    $ws1.Cells.Item(20,4).Value = 78

    $ws1.Cells.Item(21,2).Value = "06:21"
    $ws1.Cells.Item(21,3).Value = "26_HERNANDEZ"
    $ws1.Cells.Item(21,4).Value = 85

    $ws1.Cells.Item(22,1).Value = "05:23:04"
    $ws1.Cells.Item(22,2).Value = "06:24"
    $ws1.Cells.Item(22,3).Value = "11_ETCHEVERRY"
    $ws1.Cells.Item(22,4).Value = 61
    $ws1.Cells.Item(22,5).Value = "LP1912"

    $ws1.Cells.Item(23,1).Value = "05:23:04"
    $ws1.Cells.Item(23,2).Value = "06:27"
    $ws1.Cells.Item(23,3).Value = "23_HERNANDEZ"
    $ws1.Cells.Item(23,4).Value = 64
    $ws1.Cells.Item(23,5).Value = "LP1912"

    $ws1.Cells.Item(24,1).Value = "04:56:49"
    $ws1.Cells.Item(24,2).Value = "06:29"
    $ws1.Cells.Item(24,3).Value = "86_EST CHICA-ESC AGRARIA"
    $ws1.Cells.Item(24,4).Value = 93
    $ws1.Cells.Item(24,5).Value = "LP1912"

    $ws1.Cells.Item(25,1).Value = "05:23:04"
    $ws1.Cells.Item(25,2).Value = "06:31"
    $ws1.Cells.Item(25,3).Value = "17X38_ROMERO"
    $ws1.Cells.Item(25,4).Value = 68
    $ws1.Cells.Item(25,5).Value = "LP1912"

    $ws1.Cells.Item(26,1).Value = "05:23:04"
    $ws1.Cells.Item(26,2).Value = "06:31"
    $ws1.Cells.Item(26,3).Value = "16_SANTA ANA"
    $ws1.Cells.Item(26,4).Value = 68
    $ws1.Cells.Item(26,5).Value = "LP1912"

    $ws1.Cells.Item(27,1).Value = "05:23:04"
    $ws1.Cells.Item(27,2).Value = "06:39"
    $ws1.Cells.Item(27,3).Value = "225_C ROCA-H SUR"
    $ws1.Cells.Item(27,4).Value = 76
    $ws1.Cells.Item(27,5).Value = "LP1912"

    $ws1.Cells.Item(28,1).Value = "04:56:49"
    $ws1.Cells.Item(28,2).Value = "06:44"
    $ws1.Cells.Item(28,3).Value = "225_C ROCA-H SUR"
    $ws1.Cells.Item(28,4).Value = 108
    $ws1.Cells.Item(28,5).Value = "LP1912"

    $ws1.Cells.Item(29,1).Value = "04:56:49"
    $ws1.Cells.Item(29,2).Value = "06:46"
    $ws1.Cells.Item(29,3).Value = "215C_EL PATO"
    $ws1.Cells.Item(29,4).Value = 110
    $ws1.Cells.Item(29,5).Value = "LP1912"

    $ws1.Cells.Item(30,1).Value = "05:23:04"
    $ws1.Cells.Item(30,2).Value = "06:54"
    $ws1.Cells.Item(30,3).Value = "14_ABASTO"
    $ws1.Cells.Item(30,4).Value = 91
    $ws1.Cells.Item(30,5).Value = "LP1912"

    $ws1.Cells.Item(31,1).Value = "05:23:04"
    $ws1.Cells.Item(31,2).Value = "07:01"
    $ws1.Cells.Item(31,3).Value = "16_SANTA ANA"
    $ws1.Cells.Item(31,4).Value = 98
    $ws1.Cells.Item(31,5).Value = "LP1912"

    $ws1.Cells.Item(32,1).Value = "05:23:04"
    $ws1.Cells.Item(32,2).Value = "07:04"
    $ws1.Cells.Item(32,3).Value = "225_GOMEZ"
    $ws1.Cells.Item(32,4).Value = 101
    $ws1.Cells.Item(32,5).Value = "LP1912"

    $ws1.Cells.Item(33,1).Value = "05:23:04"
    $ws1.Cells.Item(33,2).Value = "07:07"
    $ws1.Cells.Item(33,3).Value = "215C_EL PATO"
    $ws1.Cells.Item(33,4).Value = 104
    $ws1.Cells.Item(33,5).Value = "LP1912"

    $ws1.Cells.Item(34,1).Value = "05:23:04"
    $ws1.Cells.Item(34,2).Value = "07:14"
    $ws1.Cells.Item(34,3).Value = "14X44_ABASTO"
    $ws1.Cells.Item(34,4).Value = 111
    $ws1.Cells.Item(34,5).Value = "LP1912"

    $ws1.Cells.Item(35,1).Value = "05:23:04"
    $ws1.Cells.Item(35,2).Value = "07:21"
    $ws1.Cells.Item(35,3).Value = "215A_EL PATO"
    $ws1.Cells.Item(35,4).Value = 118
    $ws1.Cells.Item(35,5).Value = "LP1912"

# --- Sheet "LP1912-215" ------------------------------------------------
    $ws2.Cells.Item(2,1).Value = "Última actualización: 05:23:04"

    $ws2.Cells.Item(3,1).Value = "Total filas: 6"

    $ws2.Cells.Item(8,1).Value = "05:23:04"
    $ws2.Cells.Item(8,4).Value = 48

    $ws2.Cells.Item(10,1).Value = "05:23:04"
    $ws2.Cells.Item(10,2).Value = "07:07"
    $ws2.Cells.Item(10,3).Value = "215C_EL PATO"
    $ws2.Cells.Item(10,4).Value = 104
    $ws2.Cells.Item(10,5).Value = "LP1912"

    $ws2.Cells.Item(11,1).Value = "05:23:04"
    $ws2.Cells.Item(11,2).Value = "07:21"
    $ws2.Cells.Item(11,3).Value = "215A_EL PATO"
    $ws2.Cells.Item(11,4).Value = 118
    $ws2.Cells.Item(11,5).Value = "LP1912"

# --- Sheet "6203-6173" --------------------------------------------------
    $ws3.Cells.Item(2,1).Value = "Última actualización: 05:23:04"
